# LMS-2523 Update BaSynthec Validation
# Fix the strain identifier on the "openbis-data" sheet and restore the
# active selection/cursor position that Excel recorded on last save.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("openbis-data")

# The strain id "MGP1000" was missing its "JJS-" prefix.
$wsData.Range("A3").Value = "JJS-MGP1000"

# Move/restore the active cell on the data sheet (was A6, now A16).
$wsData.Activate()
$wsData.Range("A16").Select() | Out-Null
